$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1320.8422
$ws.Range("J17").Value = 1125.1082
$ws.Range("L17").Value = 3375.3246
$ws.Range("N17").Value = -3711.3246
$ws.Range("H19").Value = 1730
$ws.Range("I19").Value = 823.4545000000001
$ws.Range("K19").Value = 823.4545000000001
$ws.Range("M19").Value = -648.4545000000001
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("M34").ClearContents()
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("M36").ClearContents()
$ws.Range("H38").Value = 391.83334
$ws.Range("I38").Value = 391.83334
$ws.Range("K38").Value = 1175.50002
$ws.Range("M38").Value = -803.5000199999999
$ws.Range("H40").Value = 1078.9445
$ws.Range("I40").Value = 1076.3103
$ws.Range("K40").Value = 1076.3103
$ws.Range("M40").Value = -901.3103000000001
$ws.Range("H51").Value = 3117.5
$ws.Range("I51").Value = 1450
$ws.Range("J51").Value = 3355.7144
$ws.Range("K51").Value = 1450
$ws.Range("L51").Value = 3355.7144
$ws.Range("M51").Value = -966
$ws.Range("N51").Value = -4323.7144
$ws.Range("H53").Value = 2579.3635
$ws.Range("I53").Value = 2424.5715
$ws.Range("J53").Value = 2850.25
$ws.Range("K53").Value = 2424.5715
$ws.Range("L53").Value = 2850.25
$ws.Range("M53").Value = -1787.5715
$ws.Range("N53").Value = -4124.25
$ws.Range("H116").Value = 16699.9
$ws.Range("I116").Value = 51000
$ws.Range("K116").Value = 51000
$ws.Range("M116").Value = -47558
$ws.Range("H127").Value = 1737.3077
$ws.Range("I127").Value = 1767.2307
$ws.Range("J127").Value = 1707.3846
$ws.Range("K127").Value = 5301.6921
$ws.Range("L127").Value = 5122.1538
$ws.Range("M127").Value = -341.6921000000002
$ws.Range("N127").Value = -15042.1538
$ws.Range("H129").Value = 915.4783
$ws.Range("J129").Value = 980.2564
$ws.Range("L129").Value = 2940.7692
$ws.Range("N129").Value = -12940.7692
$ws.Range("H131").Value = 2158.5334
$ws.Range("I131").Value = 643.0714
$ws.Range("K131").Value = 1929.2142
$ws.Range("M131").Value = 3110.7858
$ws.Range("H132").Value = 7577249
$ws.Range("I132").Value = 8773307
$ws.Range("J132").Value = 2217
$ws.Range("K132").Value = 26319921
$ws.Range("L132").Value = 6651
$ws.Range("M132").Value = -26317391
$ws.Range("N132").Value = -11711
$ws.Range("H135").Value = 479.19354
$ws.Range("I135").Value = 439.82758
$ws.Range("K135").Value = 3958.44822
$ws.Range("M135").Value = -1423.44822
$ws.Range("H137").Value = 1367.3158
$ws.Range("I137").Value = 928.53845
$ws.Range("K137").Value = 2785.61535
$ws.Range("M137").Value = -235.61535
$ws.Range("H138").Value = 1438.965
$ws.Range("I138").Value = 1364.717
$ws.Range("J138").Value = 2422.75
$ws.Range("K138").Value = 4094.151
$ws.Range("L138").Value = 7268.25
$ws.Range("M138").Value = 1045.849
$ws.Range("N138").Value = -17548.25
$ws.Range("H139").Value = 48787.375
$ws.Range("J139").Value = 48787.375
$ws.Range("L139").Value = 48787.375
$ws.Range("N139").Value = -59067.375
$ws.Range("H140").Value = 59121.285
$ws.Range("J140").Value = 59121.285
$ws.Range("L140").Value = 59121.285
$ws.Range("N140").Value = -69481.285
$ws.Range("H141").Value = 584664.6
$ws.Range("I141").Value = 718622.1
$ws.Range("J141").Value = 4182
$ws.Range("K141").Value = 2155866.3
$ws.Range("L141").Value = 12546
$ws.Range("M141").Value = -2150686.3
$ws.Range("N141").Value = -22906

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1111711
$ws.Range("I2").Value = 1389388.8
$ws.Range("K2").Value = 1389388.8
$ws.Range("M2").Value = -1389275.8
$ws.Range("H32").Value = 3155.9302
$ws.Range("I32").Value = 2502.1428
$ws.Range("J32").Value = 8749.444
$ws.Range("K32").Value = 2502.1428
$ws.Range("L32").Value = 8749.444
$ws.Range("M32").Value = -2215.1428
$ws.Range("N32").Value = -9323.444
$ws.Range("H45").Value = 1784.9286
$ws.Range("I45").Value = 1599.4
$ws.Range("K45").Value = 1599.4
$ws.Range("M45").Value = -1222.4
$ws.Range("H61").Value = 50001156
$ws.Range("I61").Value = 27779046
$ws.Range("J61").Value = 250000140
$ws.Range("K61").Value = 27779046
$ws.Range("L61").Value = 250000140
$ws.Range("M61").Value = -27778834
$ws.Range("N61").Value = -250000564
$ws.Range("H74").Value = 1271.3954
$ws.Range("I74").Value = 993.9643
$ws.Range("J74").Value = 1789.2667
$ws.Range("K74").Value = 993.9643
$ws.Range("L74").Value = 1789.2667
$ws.Range("M74").Value = -119.9643
$ws.Range("N74").Value = -3537.2667
$ws.Range("H77").Value = 1271.3954
$ws.Range("I77").Value = 993.9643
$ws.Range("J77").Value = 1789.2667
$ws.Range("K77").Value = 4969.8215
$ws.Range("L77").Value = 8946.333499999999
$ws.Range("M77").Value = -601.8215
$ws.Range("N77").Value = -17682.3335
$ws.Range("H94").Value = 99996.5
$ws.Range("J94").Value = 99996.5
$ws.Range("L94").Value = 99996.5
$ws.Range("N94").Value = -101798.5
$ws.Range("H116").Value = 1111711
$ws.Range("I116").Value = 1389388.8
$ws.Range("K116").Value = 1389388.8
$ws.Range("M116").Value = -1387094.8
$ws.Range("H132").Value = 1417.4348
$ws.Range("I132").Value = 1004.4054
$ws.Range("J132").Value = 3115.4443
$ws.Range("K132").Value = 3013.2162
$ws.Range("L132").Value = 9346.332900000001
$ws.Range("M132").Value = -483.2161999999998
$ws.Range("N132").Value = -14406.3329
$ws.Range("H136").Value = 50001156
$ws.Range("I136").Value = 27779046
$ws.Range("J136").Value = 250000140
$ws.Range("K136").Value = 83337138
$ws.Range("L136").Value = 750000420
$ws.Range("M136").Value = -83334588
$ws.Range("N136").Value = -750005520

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1111711
$ws.Range("I3").Value = 1389388.8
$ws.Range("K3").Value = 1389388.8
$ws.Range("M3").Value = -1389274.8
$ws.Range("H80").Value = 10510.6
$ws.Range("J80").Value = 14825.857
$ws.Range("L80").Value = 14825.857
$ws.Range("N80").Value = -16821.857
$ws.Range("H83").Value = 10510.6
$ws.Range("J83").Value = 14825.857
$ws.Range("L83").Value = 74129.285
$ws.Range("N83").Value = -84113.285
$ws.Range("H107").Value = 1764.3334
$ws.Range("I107").Value = 1821.75
$ws.Range("J107").Value = 1649.5
$ws.Range("K107").Value = 1821.75
$ws.Range("L107").Value = 1649.5
$ws.Range("M107").Value = 98.25
$ws.Range("N107").Value = -5489.5
$ws.Range("H130").Value = 32499.25
$ws.Range("J130").Value = 32499.25
$ws.Range("L130").Value = 32499.25
$ws.Range("N130").Value = -42539.25
$ws.Range("H134").Value = 4743.1567
$ws.Range("I134").Value = 3959.5957
$ws.Range("K134").Value = 11878.7871
$ws.Range("M134").Value = -9343.7871

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2142.2856
$ws.Range("I16").Value = 799
$ws.Range("J16").Value = 3149.75
$ws.Range("K16").Value = 799
$ws.Range("L16").Value = 3149.75
$ws.Range("M16").Value = -512
$ws.Range("N16").Value = -3723.75
$ws.Range("H31").Value = 1985722
$ws.Range("I31").Value = 4202858.5
$ws.Range("J31").Value = 1968.3158
$ws.Range("K31").Value = 4202858.5
$ws.Range("L31").Value = 1968.3158
$ws.Range("M31").Value = -4202563.5
$ws.Range("N31").Value = -2558.3158
$ws.Range("H34").Value = 1985722
$ws.Range("I34").Value = 4202858.5
$ws.Range("J34").Value = 1968.3158
$ws.Range("K34").Value = 4202858.5
$ws.Range("L34").Value = 1968.3158
$ws.Range("M34").Value = -4202656.5
$ws.Range("N34").Value = -2372.3158
$ws.Range("H58").Value = 3485.3845
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 3485.3845
$ws.Range("K58").Value = 0
$ws.Range("M58").ClearContents()
$ws.Range("N58").Value = -3891.3845
$ws.Range("H88").Value = 59000
$ws.Range("J88").Value = 59000
$ws.Range("L88").Value = 59000
$ws.Range("N88").Value = -59812
$ws.Range("H91").Value = 59000
$ws.Range("J91").Value = 59000
$ws.Range("L91").Value = 59000
$ws.Range("N91").Value = -61808
$ws.Range("H113").Value = 2142.2856
$ws.Range("I113").Value = 799
$ws.Range("J113").Value = 3149.75
$ws.Range("K113").Value = 799
$ws.Range("L113").Value = 3149.75
$ws.Range("M113").Value = 1371
$ws.Range("N113").Value = -7489.75
$ws.Range("H132").Value = 1779.1945
$ws.Range("I132").Value = 1038.3103
$ws.Range("K132").Value = 3114.9309
$ws.Range("M132").Value = -584.9309000000003
$ws.Range("H134").Value = 1784.8823
$ws.Range("I134").Value = 1521.5
$ws.Range("J134").Value = 5999
$ws.Range("K134").Value = 4564.5
$ws.Range("L134").Value = 17997
$ws.Range("M134").Value = -2029.5
$ws.Range("N134").Value = -23067
$ws.Range("H136").Value = 3485.3845
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 3485.3845
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()
$ws.Range("N136").Value = -15556.1535

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 67.89474
$ws.Range("J2").Value = 33
$ws.Range("L2").Value = 198
$ws.Range("N2").Value = -424
$ws.Range("H4").Value = 256767.27
$ws.Range("I4").Value = 6771.933
$ws.Range("K4").Value = 20315.799
$ws.Range("M4").Value = -20203.799
$ws.Range("H20").Value = 1000
$ws.Range("J20").Value = 1000
$ws.Range("L20").Value = 3000
$ws.Range("N20").Value = -3454
$ws.Range("H22").Value = 667000.3
$ws.Range("I22").Value = 500500.5
$ws.Range("J22").Value = 1000000
$ws.Range("K22").Value = 1501501.5
$ws.Range("L22").Value = 3000000
$ws.Range("M22").Value = -1501332.5
$ws.Range("N22").Value = -3000338
$ws.Range("H23").Value = 100
$ws.Range("I23").Value = 100
$ws.Range("J23").Value = 100
$ws.Range("K23").Value = 300
$ws.Range("L23").Value = 300
$ws.Range("M23").Value = -65
$ws.Range("N23").Value = -770
$ws.Range("H27").Value = 667000.3
$ws.Range("I27").Value = 500500.5
$ws.Range("J27").Value = 1000000
$ws.Range("K27").Value = 1501501.5
$ws.Range("L27").Value = 3000000
$ws.Range("M27").Value = -1501399.5
$ws.Range("N27").Value = -3000204
$ws.Range("H32").Value = 1000
$ws.Range("J32").Value = 1000
$ws.Range("L32").Value = 3000
$ws.Range("N32").Value = -3566
$ws.Range("H37").Value = 100000
$ws.Range("J37").Value = 100000
$ws.Range("L37").Value = 300000
$ws.Range("N37").Value = -300224
$ws.Range("H38").Value = 372.875
$ws.Range("I38").Value = 176.2
$ws.Range("J38").Value = 700.6667
$ws.Range("K38").Value = 528.5999999999999
$ws.Range("L38").Value = 2102.0001
$ws.Range("M38").Value = -181.5999999999999
$ws.Range("N38").Value = -2796.0001
$ws.Range("H44").Value = 1081.6666
$ws.Range("I44").Value = 500
$ws.Range("J44").Value = 1372.5
$ws.Range("K44").Value = 1500
$ws.Range("L44").Value = 4117.5
$ws.Range("M44").Value = -1102
$ws.Range("N44").Value = -4913.5
$ws.Range("H46").Value = 800
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()
$ws.Range("H62").Value = 4000
$ws.Range("J62").Value = 4000
$ws.Range("L62").Value = 12000
$ws.Range("N62").Value = -13372
$ws.Range("H65").Value = 4000
$ws.Range("J65").Value = 4000
$ws.Range("L65").Value = 36000
$ws.Range("N65").Value = -42864
$ws.Range("H118").Value = 1739.9
$ws.Range("I118").Value = 1523.2
$ws.Range("K118").Value = 4569.6
$ws.Range("M118").Value = -3326.6
$ws.Range("H122").Value = 903.82855
$ws.Range("J122").Value = 922.5357
$ws.Range("L122").Value = 8302.8213
$ws.Range("N122").Value = -13202.8213
$ws.Range("H129").Value = 80996.336
$ws.Range("J129").Value = 145393.6
$ws.Range("L129").Value = 436180.8
$ws.Range("N129").Value = -446180.8
$ws.Range("H131").Value = 8186.6313
$ws.Range("J131").Value = 9419.134
$ws.Range("L131").Value = 28257.402
$ws.Range("N131").Value = -38337.402
$ws.Range("H137").Value = 4722.4
$ws.Range("I137").Value = 2036.6666
$ws.Range("J137").Value = 5393.8335
$ws.Range("K137").Value = 6109.9998
$ws.Range("L137").Value = 16181.5005
$ws.Range("M137").Value = -1009.9998
$ws.Range("N137").Value = -26381.5005
$ws.Range("H139").Value = 8766.177
$ws.Range("I139").Value = 12193.182
$ws.Range("K139").Value = 36579.546
$ws.Range("M139").Value = -31439.546

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 47.958332
$ws.Range("I2").Value = 10.076923
$ws.Range("J2").Value = 92.72727
$ws.Range("K2").Value = 10.076923
$ws.Range("L2").Value = 92.72727
$ws.Range("M2").Value = 102.923077
$ws.Range("N2").Value = -318.72727
$ws.Range("H11").Value = 4524073.5
$ws.Range("I11").Value = 5276363.5
$ws.Range("J11").Value = 1420877.2
$ws.Range("K11").Value = 5276363.5
$ws.Range("L11").Value = 1420877.2
$ws.Range("M11").Value = -5276224.5
$ws.Range("N11").Value = -1421155.2
$ws.Range("H70").Value = 19361.2
$ws.Range("I70").Value = 41953
$ws.Range("K70").Value = 41953
$ws.Range("M70").Value = -41683
$ws.Range("H73").Value = 19361.2
$ws.Range("I73").Value = 41953
$ws.Range("K73").Value = 41953
$ws.Range("M73").Value = -41017
$ws.Range("H102").Value = 1660.8928
$ws.Range("I102").Value = 1615.0385
$ws.Range("K102").Value = 1615.0385
$ws.Range("M102").Value = 6.961499999999887
$ws.Range("H122").Value = 1714.7
$ws.Range("I122").Value = 1711.625
$ws.Range("J122").Value = 1727
$ws.Range("K122").Value = 5134.875
$ws.Range("L122").Value = 5181
$ws.Range("M122").Value = -2684.875
$ws.Range("N122").Value = -10081
$ws.Range("H132").Value = 1014624.06

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 392500
$ws.Range("J2").Value = 70000
$ws.Range("L2").Value = 70000
$ws.Range("N2").Value = -70224
$ws.Range("H7").Value = 3454.077
$ws.Range("I7").Value = 2414.8572
$ws.Range("K7").Value = 2414.8572
$ws.Range("M7").Value = -2302.8572
$ws.Range("H22").Value = 1909.2142
$ws.Range("I22").Value = 2254.8333
$ws.Range("J22").Value = 1650
$ws.Range("K22").Value = 2254.8333
$ws.Range("L22").Value = 1650
$ws.Range("M22").Value = -1959.8333
$ws.Range("N22").Value = -2240
$ws.Range("H27").Value = 1909.2142
$ws.Range("I27").Value = 2254.8333
$ws.Range("J27").Value = 1650
$ws.Range("K27").Value = 2254.8333
$ws.Range("L27").Value = 1650
$ws.Range("M27").Value = -2147.8333
$ws.Range("N27").Value = -1864
$ws.Range("H29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("N29").ClearContents()
$ws.Range("H44").Value = 20000
$ws.Range("J44").Value = 20000
$ws.Range("L44").Value = 20000
$ws.Range("N44").Value = -20912
$ws.Range("H45").Value = 40997.5
$ws.Range("I45").Value = 40997.5
$ws.Range("K45").Value = 40997.5
$ws.Range("M45").Value = -40590.5
$ws.Range("H55").Value = 418.1579
$ws.Range("I55").Value = 377.27274
$ws.Range("K55").Value = 377.27274
$ws.Range("M55").Value = -204.27274
$ws.Range("H93").Value = 1073.2667
$ws.Range("I93").Value = 741.3
$ws.Range("J93").Value = 1737.2
$ws.Range("K93").Value = 741.3
$ws.Range("L93").Value = 1737.2
$ws.Range("M93").Value = 506.7
$ws.Range("N93").Value = -4233.2
$ws.Range("H108").Value = 68000
$ws.Range("J108").Value = 68000
$ws.Range("L108").Value = 68000
$ws.Range("N108").Value = -75680
$ws.Range("H110").Value = 6999
$ws.Range("J110").Value = 6999
$ws.Range("L110").Value = 6999
$ws.Range("N110").Value = -15179
$ws.Range("H122").Value = 6168.1
$ws.Range("I122").Value = 10625.667
$ws.Range("K122").Value = 31877.001
$ws.Range("M122").Value = -29427.001
$ws.Range("H126").Value = 3454.077
$ws.Range("I126").Value = 2414.8572
$ws.Range("K126").Value = 7244.571599999999
$ws.Range("M126").Value = -4774.571599999999
$ws.Range("H132").Value = 1279.6934
$ws.Range("I132").Value = 907.98114
$ws.Range("K132").Value = 2723.94342
$ws.Range("M132").Value = -193.9434200000001
$ws.Range("H136").Value = 1696.678
$ws.Range("I136").Value = 1021.5208
$ws.Range("K136").Value = 3064.5624
$ws.Range("M136").Value = -514.5623999999998

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 250052500
$ws.Range("J2").Value = 70003
$ws.Range("L2").Value = 70003
$ws.Range("N2").Value = -70227
$ws.Range("H4").Value = 6000
$ws.Range("J4").Value = 6000
$ws.Range("L4").Value = 6000
$ws.Range("N4").Value = -6226
$ws.Range("H5").Value = 8000
$ws.Range("I5").Value = 8000
$ws.Range("J5").Value = 8000
$ws.Range("K5").Value = 8000
$ws.Range("L5").Value = 8000
$ws.Range("M5").Value = -7888
$ws.Range("N5").Value = -8224
$ws.Range("H18").Value = 12856.429
$ws.Range("I18").Value = 12000
$ws.Range("J18").Value = 12999.167
$ws.Range("K18").Value = 12000
$ws.Range("L18").Value = 12999.167
$ws.Range("M18").Value = -11827
$ws.Range("N18").Value = -13345.167
$ws.Range("H81").Value = 783.3333
$ws.Range("I81").Value = 300
$ws.Range("J81").Value = 1025
$ws.Range("K81").Value = 600
$ws.Range("L81").Value = 2050
$ws.Range("M81").Value = 461
$ws.Range("N81").Value = -4172
$ws.Range("H84").Value = 783.3333
$ws.Range("I84").Value = 300
$ws.Range("J84").Value = 1025
$ws.Range("K84").Value = 3000
$ws.Range("L84").Value = 10250
$ws.Range("M84").Value = 2304
$ws.Range("N84").Value = -20858
$ws.Range("H100").Value = 1016.5
$ws.Range("I100").Value = 993.7273
$ws.Range("J100").Value = 1100
$ws.Range("K100").Value = 1987.4546
$ws.Range("L100").Value = 2200
$ws.Range("M100").Value = -1446.4546
$ws.Range("N100").Value = -3282
$ws.Range("H107").Value = 853.3333
$ws.Range("I107").Value = 315.5
$ws.Range("J107").Value = 1391.1666
$ws.Range("K107").Value = 946.5
$ws.Range("L107").Value = 4173.4998
$ws.Range("M107").Value = 973.5
$ws.Range("N107").Value = -8013.4998
$ws.Range("H132").Value = 1187.0759
$ws.Range("I132").Value = 775.9365
$ws.Range("J132").Value = 2805.9375
$ws.Range("K132").Value = 2327.8095
$ws.Range("L132").Value = 8417.8125
$ws.Range("M132").Value = 202.1904999999997
$ws.Range("N132").Value = -13477.8125
$ws.Range("H136").Value = 11112737
$ws.Range("I136").Value = 13551650
$ws.Range("J136").Value = 2134.4443
$ws.Range("K136").Value = 40654950
$ws.Range("L136").Value = 6403.3329
$ws.Range("M136").Value = -40652400
$ws.Range("N136").Value = -11503.3329
